# Apply updated crypto price/volume figures to sheet1 (as described by the commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimals with a trailing zero (e.g. "226.10").
# Excel auto-converts such strings to numbers on assignment and would silently
# drop the trailing zero (General format), so force those specific cells to text
# first, then restore the default (Normal) style so no stray formatting remains.
$forceTextCells = @("D19", "D28", "D33", "D36", "D47", "D51")
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.774.64"
$ws.Range("E2").Value = "  +6.42%  "
$ws.Range("D3").Value = "1.737.23"
$ws.Range("E3").Value = "  +5.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "227.28"
$ws.Range("E5").Value = "  +3.97%  "
$ws.Range("D6").Value = "0.5457"
$ws.Range("E6").Value = "  +3.89%  "
$ws.Range("D8").Value = "0.2757"
$ws.Range("E8").Value = "  +3.05%  "
$ws.Range("D9").Value = "0.06725"
$ws.Range("E9").Value = "  +5.60%  "
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("D11").Value = "0.07777"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "1.737.47"
$ws.Range("E13").Value = "  +4.92%  "
$ws.Range("D14").Value = "1.976.23"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("D15").Value = "0.5975"
$ws.Range("E15").Value = "  +6.29%  "
$ws.Range("D16").Value = "0.0₅8423"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "69.16"
$ws.Range("E17").Value = "  +5.63%  "
$ws.Range("D18").Value = "27.783.06"
$ws.Range("E18").Value = "  +6.48%  "
$ws.Range("D19").Value = "226.10"
$ws.Range("E19").Value = "  +18.32%  "
$ws.Range("D20").Value = "4.826"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("E22").Value = "  +5.24%  "
$ws.Range("D23").Value = "6.222"
$ws.Range("E23").Value = "  +4.19%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "146.77"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "0.1247"
$ws.Range("E26").Value = "  +3.80%  "
$ws.Range("D27").Value = "1.701"
$ws.Range("E27").Value = "  +13.58%  "
$ws.Range("D28").Value = "7.450"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("D29").Value = "17.16"
$ws.Range("E29").Value = "  +7.60%  "
$ws.Range("D30").Value = "0.05655"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("D32").Value = "3.686"
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("D33").Value = "3.510"
$ws.Range("E33").Value = "  +3.79%  "
$ws.Range("D34").Value = "1.674"
$ws.Range("E34").Value = "  +6.04%  "
$ws.Range("D35").Value = "0.9761"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "2.860"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "2.448"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").Value = "0.5943"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("D39").Value = "0.01667"
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("D40").Value = "5.873"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "0.8475"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "1.047.64"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "101.94"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "1.881.44"
$ws.Range("E45").Value = "  +5.08%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +12.52%  "
$ws.Range("D47").Value = "59.20"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("D48").Value = "8.255"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").Value = "0.4437"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").Value = "0.9987"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "0.05310"
$ws.Range("E51").Value = "  -0.69%  "

foreach ($ref in $forceTextCells) {
    $ws.Range($ref).Style = "Normal"
}
